# "Concept 1 for report" -- adds a "Part"/structure column to the Quotes
# sheet, fixes the Marvin Minsky quote punctuation, and adds a new
# Yuval Noah Harrari / Homo Deus quote row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quotes")

# Insert a new first column (Report/Part) in front of the existing
# Auteur/Quote/Quote english/Bron columns, shifting everything right.
$ws.Columns.Item(1).Insert()

# Give the new column a width roughly matching the other label columns
# (closest achievable value via ColumnWidth rounding).
$ws.Columns.Item(1).ColumnWidth = 16.86

# Part labels for the existing quote rows (filled in roughly the order
# the parts were drafted).
$ws.Range("A2").Value = "Part I"
$ws.Range("A4").Value = "Voorwoord"
$ws.Range("A5").Value = "Voorwoord"
$ws.Range("A7").Value = "Part III"
$ws.Range("A8").Value = "Part IV"
$ws.Range("A3").Value = "Part II"

# Header, added last.
$ws.Range("A1").Value = "Report"

# Fix the Marvin Minsky quote text (em dash -> period, capitalised "Like").
$ws.Range("D3").Value = "We rarely recognize how wonderful it is that a person can traverse an entire lifetime without making a single really serious mistake. Like putting a fork in one's eye or using a window instead of a door."

# New quote row: Yuval Noah Harrari / Homo Deus.
$ws.Range("B11").Value = "Yuval Noah Harrari"
$ws.Range("E11").Value = "Homo Deus"
$ws.Range("D11").Value = "Studying history will not tell us what to choose, but at least gives us more options to choose from."
$ws.Range("A11").Value = "Part V"

# Move the selection, matching the saved file's cursor position.
$ws.Range("A12").Select()
